$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "PDL Date" column (B) for every data row (2 through 407) changes
# from "04012025" to "99999999". These values are dates stored as plain
# text (e.g. inline/string cells), not real numbers, so a leading
# apostrophe is used to force text entry instead of letting Excel
# auto-convert the digit string into a numeric value.
$ws.Range("B2:B407").Value = "'99999999"
